# Weekly update of "Perejil" (parsley) price records for
# Terminal Hortofrutícola Agro Chillán.
#
# Two brand-new weekly rows are prepended to the data block (new rows
# 107-108), and every existing data row (old rows 107-128) shifts down by
# two positions (new rows 109-130). The per-sheet constant columns
# (Mercado ID, Mercado, Región, Codreg, Categoría ID, Categoría, Variedad,
# Unidad de comercialización, Kg o Unidades, Clasificación) are identical
# for every row in this sheet, so the two new rows reuse them verbatim.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Make room for the two new weekly records: push the existing data
# (old rows 107:128) down to new rows 109:130.
$ws.Range("A107:A108").EntireRow.Insert()

# column layout:
#  1 Mercado ID          7  Categoría           13 Precio promedio ponderado
#  2 Mercado             8  Variedad            14 Unidad de comercialización
#  3 Región              9  Calidad             15 Origen
#  4 Fecha              10  Volumen             16 Precio $/Kg
#  5 Codreg             11  Precio mínimo       17 Kg o Unidades
#  6 Categoría ID       12  Precio máximo       18 Clasificación

$rows = @(
  @(107, 45211, "Primera", 120, 1300, 1500, 1400, "Región de Ñuble", 1400),
  @(108, 45211, "Segunda", 200, 1000, 1000, 1000, "Región de Ñuble", 1000)
)

foreach ($r in $rows) {
    $rowNum = $r[0]

    $ws.Cells.Item($rowNum, 1).Value  = 7
    $ws.Cells.Item($rowNum, 2).Value  = "Terminal Hortofrutícola Agro Chillán"
    $ws.Cells.Item($rowNum, 3).Value  = "Ñuble"
    $ws.Cells.Item($rowNum, 4).Value  = $r[1]
    $ws.Cells.Item($rowNum, 5).Value  = 16
    $ws.Cells.Item($rowNum, 6).Value  = 100112044
    $ws.Cells.Item($rowNum, 7).Value  = "Perejil"
    $ws.Cells.Item($rowNum, 8).Value  = "Sin especificar"
    $ws.Cells.Item($rowNum, 9).Value  = $r[2]
    $ws.Cells.Item($rowNum, 10).Value = $r[3]
    $ws.Cells.Item($rowNum, 11).Value = $r[4]
    $ws.Cells.Item($rowNum, 12).Value = $r[5]
    $ws.Cells.Item($rowNum, 13).Value = $r[6]
    $ws.Cells.Item($rowNum, 14).Value = "`$/atado 0,5 a 1 kilo"
    $ws.Cells.Item($rowNum, 15).Value = $r[7]
    $ws.Cells.Item($rowNum, 16).Value = $r[8]
    $ws.Cells.Item($rowNum, 17).Value = 1
    $ws.Cells.Item($rowNum, 18).Value = "Hortaliza"
}
